$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FolioChildLocations")
$ws.Activate()

$ws.Columns("R").Insert()
$ws.Columns("R").ColumnWidth = 16.81640625

$ws.Range("R1").Value = "ADDITIONAL_DETAILS_OF_TRANSITION_PLAN"
$ws.Range("R2").Value = "n/a"
$ws.Range("R3").Value = "n/a"
$ws.Range("R4").Value = "n/a"
$ws.Range("R5").Value = "n/a"

$ws.Range("R2:R5").Select() | Out-Null

Write-Host "done"
